# daily auto push: 2026-01-09 02:27 UTC
# A new sensor reading for 2026/01/09 (time 8, ranking 182) was appended
# right after the existing 2026/01/09 row (old row 611), pushing every
# subsequent row down by one (old row 611 -> new row 612, ..., old row
# 652 -> new row 653).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 611; this shifts rows 611..652 down to
# 612..653 and keeps all of their original formatting/values intact.
$ws.Rows("611:611").Insert()

# Fill the newly inserted row with the new observation. The date/weekday
# columns in this sheet store plain text (e.g. "2026/01/09"), not real
# Excel dates, so a leading apostrophe forces text entry instead of
# Excel's automatic date parsing; ClearFormats() then drops the
# "entered via quote prefix" formatting flag so the cell's style matches
# its neighbours (no explicit style id), exactly like the rest of the
# column.
$ws.Range("A611").Value = "'2026/01/09"
$ws.Range("B611").Value = "金"
$ws.Range("C611").Value = 8
$ws.Range("D611").Value = 182
$ws.Range("A611").ClearFormats()
